# Generate Report for Archive
# The localization run moved from "Ready for handoff" to "In Translation",
# so every cell carrying that status text is updated, and the Status /
# language columns that used to size themselves to the old (longer) text
# are re-sized to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: the zh-cn / de-de columns hold the status text directly ---
$wsOverview = $wb.Worksheets.Item("Overview")
$ovRows = $wsOverview.UsedRange.Rows.Count
for ($r = 2; $r -le $ovRows; $r++) {
    foreach ($colLetter in @("E", "F")) {
        $cell = $wsOverview.Range($colLetter + $r)
        if ($cell.Value() -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}
$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null
# Match the precise content-fit width Excel computes for the new text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.576851254417766
$wsOverview.Columns.Item(6).ColumnWidth = 12.576851254417766

# --- Per-locale sheets: column C ("Status") holds the status text ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $rows; $r++) {
        $cell = $ws.Range("C" + $r)
        if ($cell.Value() -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
    $ws.Columns.Item(3).AutoFit() | Out-Null
    $ws.Columns.Item(3).ColumnWidth = 12.576851254417766
}
